$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("E42:G48").Value = "andet"
$ws.Range("E49:F49").Value = "andet"
$ws.Range("G49").Value = "andet/4"
$ws.Range("E50:F50").Value = "andet/4"
$ws.Range("I54").Value = "andet"
$ws.Range("I55").Value = "andet/2"
$ws.Range("I56").Value = "andet/4"

$ws.Range("I57").Select()
